$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new daily price record was inserted at row 124 ("Hortaliza, Feria Lagunitas
# de Puerto Montt - Zapallo italiano"), pushing the existing rows 124-221 down
# to 125-222 and extending the sheet's used range to A1:R222.
$ws.Rows.Item(124).Insert()

$ws.Cells.Item(124, 1).Value = 4
$ws.Cells.Item(124, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(124, 3).Value = "Los Lagos"
$ws.Cells.Item(124, 4).Value = 44651
$ws.Cells.Item(124, 5).Value = 10
$ws.Cells.Item(124, 6).Value = 100112032
$ws.Cells.Item(124, 7).Value = "Zapallo italiano"
$ws.Cells.Item(124, 8).Value = "Sin especificar"
$ws.Cells.Item(124, 9).Value = "Primera"
$ws.Cells.Item(124, 10).Value = 70
$ws.Cells.Item(124, 11).Value = 13000
$ws.Cells.Item(124, 12).Value = 13000
$ws.Cells.Item(124, 13).Value = 13000
$ws.Cells.Item(124, 14).Value = '$/caja 50 unidades'
$ws.Cells.Item(124, 15).Value = "Región Metropolitana"
$ws.Cells.Item(124, 16).Value = 260
$ws.Cells.Item(124, 17).Value = 50
$ws.Cells.Item(124, 18).Value = "Hortaliza"
